# Generate Report for Handoff
# Adds a new localization-status row (for file
# 532104a5-dc69-499e-9a93-80bffce03326.md) to the Overview / zh-cn / de-de
# sheets, mirroring the existing 274e3d81-... row.

$wb = $excel.ActiveWorkbook

$newGuid = "532104a5-dc69-499e-9a93-80bffce03326"
$newMd = "$newGuid.md"
$newMdDisplay = "e2e\$newGuid.md"
$commit = "fd38f606b6a1fc0524797a448be83d62db68f57f"
$newMdUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/$commit/e2e/$newGuid.md"

# ---------------------------------------------------------------------------
# Overview sheet (row 3)
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$loOverview = $wsOverview.ListObjects.Item(1)
$loOverview.ListRows.Add() | Out-Null

$wsOverview.Range("A3").Value = $newMd
$wsOverview.Range("B3").Value = $newMdDisplay
$wsOverview.Range("C3").Value = ".md"
$wsOverview.Range("D3").Value = ""
$wsOverview.Range("E3").Value = "Ready for handoff"
$wsOverview.Range("F3").Value = "Ready for handoff"
$wsOverview.Range("G3").Value = "2016-08-29 10:39:57"
$wsOverview.Range("G3").NumberFormat = "yyyy-mm-dd HH:mm:ss"

$wsOverview.Hyperlinks.Add($wsOverview.Range("B3"), $newMdUrl, "", "", $newMdDisplay) | Out-Null

# ---------------------------------------------------------------------------
# zh-cn sheet (row 3)
# ---------------------------------------------------------------------------
$zhXlf = "$newGuid.f750fb45b29ac2415fb0e572f4e51e5c55565cf3.zh-cn.xlf"

$wsZh = $wb.Worksheets.Item("zh-cn")
$loZh = $wsZh.ListObjects.Item(1)
$loZh.ListRows.Add() | Out-Null

$wsZh.Range("A3").Value = $newMd
$wsZh.Range("B3").Value = ".md"
$wsZh.Range("C3").Value = "Ready for handoff"
$wsZh.Range("D3").Value = "e2e"
$wsZh.Range("E3").Value = "ht"
$wsZh.Range("F3").Value = "'False"
$wsZh.Range("G3").Value = $zhXlf
$wsZh.Range("H3").Value = "2016-08-29 10:39:52"
$wsZh.Range("H3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZh.Range("I3").Value = ""
$wsZh.Range("J3").Value = ""
$wsZh.Range("K3").Value = "0001-01-01 00:00:00"
$wsZh.Range("K3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZh.Range("L3").Value = ""
$wsZh.Range("M3").Value = "'True"
$wsZh.Range("N3").Value = ""
$wsZh.Range("O3").Value = "'False"
$wsZh.Range("P3").Value = ""

$wsZh.Hyperlinks.Add($wsZh.Range("A3"), $newMdUrl, "", "", $newMd) | Out-Null
$wsZh.Range("A3").Style = "HyperLink"

# ---------------------------------------------------------------------------
# de-de sheet (row 3)
# ---------------------------------------------------------------------------
$deXlf = "$newGuid.f750fb45b29ac2415fb0e572f4e51e5c55565cf3.de-de.xlf"

$wsDe = $wb.Worksheets.Item("de-de")
$loDe = $wsDe.ListObjects.Item(1)
$loDe.ListRows.Add() | Out-Null

$wsDe.Range("A3").Value = $newMd
$wsDe.Range("B3").Value = ".md"
$wsDe.Range("C3").Value = "Ready for handoff"
$wsDe.Range("D3").Value = "e2e"
$wsDe.Range("E3").Value = "ht"
$wsDe.Range("F3").Value = "'False"
$wsDe.Range("G3").Value = $deXlf
$wsDe.Range("H3").Value = "2016-08-29 10:39:57"
$wsDe.Range("H3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDe.Range("I3").Value = ""
$wsDe.Range("J3").Value = ""
$wsDe.Range("K3").Value = "0001-01-01 00:00:00"
$wsDe.Range("K3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDe.Range("L3").Value = ""
$wsDe.Range("M3").Value = "'True"
$wsDe.Range("N3").Value = ""
$wsDe.Range("O3").Value = "'False"
$wsDe.Range("P3").Value = ""

$wsDe.Hyperlinks.Add($wsDe.Range("A3"), $newMdUrl, "", "", $newMd) | Out-Null
$wsDe.Range("A3").Style = "HyperLink"
